$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Stash the existing "thin left border" style (currently s="1", on
#    E1/H1/H2) onto the new cells that need it (G1/G2 "patch z" header
#    + value, O1/O2 "(stride z)" header + value) BEFORE we overwrite
#    the source cells' own contents/formats. PasteSpecial(formats)
#    reuses the existing style index instead of minting a new one.
# ------------------------------------------------------------------
$ws.Range("E1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("E1").Copy()
$ws.Range("G2").PasteSpecial(-4122)
$ws.Range("E1").Copy()
$ws.Range("O1").PasteSpecial(-4122)
$ws.Range("E1").Copy()
$ws.Range("O2").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 2) Write cell text/values. Shared-string order in the saved file is
#    driven by first-write order for genuinely new text, while text
#    that already existed keeps its relative order automatically - so
#    write the handful of still-existing labels first (in their old
#    relative order) and then the brand-new labels in the exact order
#    they need to land in the rebuilt table.
# ------------------------------------------------------------------
$ws.Range("A1").Value = "run"
$ws.Range("B1").Value = "VRAM usage GB"
$ws.Range("G1").Value = "patch z"
$ws.Range("H1").Value = "patch y"
$ws.Range("I1").Value = "patch x"
$ws.Range("R1").Value = "(resolution)"
$ws.Range("A2").Value = "230818-0?"

$ws.Range("C1").Value = "n images"
$ws.Range("E1").Value = "n patches"
$ws.Range("D1").Value = "VRAM/image"
$ws.Range("F1").Value = "VRAM/patch"
$ws.Range("O1").Value = "(stride z)"
$ws.Range("P1").Value = "(stride y)"
$ws.Range("Q1").Value = "(stride x)"
$ws.Range("K1").Value = "n raw channels"
$ws.Range("J1").Value = "pixels/patch/channel"
$ws.Range("L1").Value = "pixels/specimen_raw"

# Row 2 values / formulas
$ws.Range("B2").Value = 33
$ws.Range("C2").Value = 5
$ws.Range("D2").Formula = "=B2/C2"
$ws.Range("E2").Value = 5
$ws.Range("F2").Formula = "=B2/E2"
$ws.Range("G2").Value = 105
$ws.Range("H2").Value = 1140
$ws.Range("I2").Value = 394
$ws.Range("J2").Formula = "=G2*H2*I2"
$ws.Range("K2").Value = 3
$ws.Range("L2").Formula = "=J2*K2"
$ws.Range("O2").Value = 10
$ws.Range("P2").Value = 10
$ws.Range("Q2").Value = 10

# ------------------------------------------------------------------
# 3) Strip the leftover formatting (old "n train images" border style,
#    old "stride" header/value border styles) from cells that must be
#    plain in the new layout, by pasting in the plain (unstyled) format
#    from a cell that has always been style-less (B1). Multi-area paste
#    targets only honour the first area in this host, so do them one
#    at a time.
# ------------------------------------------------------------------
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("B1").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("B1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("B1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("B1").Copy()
$ws.Range("H2").PasteSpecial(-4122)
$ws.Range("B1").Copy()
$ws.Range("I2").PasteSpecial(-4122)
$ws.Range("B1").Copy()
$ws.Range("J2").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

# ------------------------------------------------------------------
# 4) Misc view bits that came along with the edit.
# ------------------------------------------------------------------
$null = $ws.Range("F9").Select()
